$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 84 is a new row appended to the data series (EQUIPMENT eval quarter).
# Column A holds the period-end date serial, formatted like the rest of the
# column (copy formatting from the row above so no new style is introduced);
# column B holds the plain numeric value.
$ws.Cells.Item(83, 1).Copy()
$ws.Cells.Item(84, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(84, 1).Value = 45884
$ws.Cells.Item(84, 2).Value = -0.2551464291630765
